$wb = $excel.ActiveWorkbook

# Use an already-styled cell (style index 1: bold, thin border, center/top align)
# from the existing VIF sheet as a formatting source so we don't introduce any
# new (unused) cell-format entries into styles.xml.
$ws2 = $wb.Worksheets.Item("VIF")

# Add the new worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "DTR-Features"

# Header row.
$newSheet.Range("B1").Value = "Features"
$newSheet.Range("C1").Value = "Coefficients"

# Data rows: A = index, B = feature name (existing shared strings), C = coefficient
# stored as text (shared string), matching the source workbook.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Petrol_tax"
$newSheet.Range("C2").Formula = "=""0.384"""

$newSheet.Range("A3").Value = 3
$newSheet.Range("B3").Value = "Population_Driver_licence(%)"
$newSheet.Range("C3").Formula = "=""0.303"""

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "Paved_Highways"
$newSheet.Range("C4").Formula = "=""0.163"""

$newSheet.Range("A5").Value = 1
$newSheet.Range("B5").Value = "Average_income"
$newSheet.Range("C5").Formula = "=""0.15"""

# Convert the formulas in column C into plain text values (keeps them as
# shared strings without leaving any number-format residue behind).
$newSheet.Range("C2:C5").Copy()
$newSheet.Range("C2:C5").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Apply the header + index-column formatting by copying it from the matching,
# already-present style on the VIF sheet (style index 1).
$ws2.Range("B1:C1").Copy()
$newSheet.Range("B1:C1").PasteSpecial(-4122)

$ws2.Range("A2:A5").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Keep the originally active sheet selected, leaving the workbook view as it
# was prior to this edit.
$wb.Worksheets.Item(1).Activate()
